# add function to modify cell color
#
# Reusable helper that sets a PowerPoint table cell's background to a
# solid fill color. Row/Column are 1-based, matching the native
# Table.Cell(Row, Column) indexing used by the PowerPoint object model.
# Rgb follows the classic VBA RGB() encoding (R + G*256 + B*65536), so
# pure red is 255 (0x0000FF).
function Set-TableCellColor {
    param($Table, $Row, $Column, $Rgb)

    $cell = $Table.Cell($Row, $Column)
    $cell.Shape.Fill.Solid()
    $cell.Shape.Fill.ForeColor.RGB = $Rgb
}

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$shp = $s.Shapes.Item(1)
$tbl = $shp.Table

$red = 255

Set-TableCellColor $tbl 2 2 $red
Set-TableCellColor $tbl 3 3 $red
Set-TableCellColor $tbl 4 4 $red
Set-TableCellColor $tbl 4 5 $red
Set-TableCellColor $tbl 4 6 $red
